$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 held the electrolytic-capacitor line item (C487393 / VKME1651J331MV / Ymin ...).
# The iso output was switched to ceramic caps, so that row's data is removed
# (the row itself stays in place, just emptied - rows below are untouched).
$ws.Range("A8:L8").ClearContents()

# Selection/view bookkeeping left behind by the edit in Excel: row 8 was
# selected (entire row) right before the clear, then the view scrolled back
# to the top-left corner.
$ws.Range("A1").Select()
$ws.Rows("8:8").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Activate()
